$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at H (existing H..N shift to I..O, P..Q shift to Q..R)
$ws.Columns.Item(8).Insert()

# Give the new column H the same width as its neighboring E:G columns
$ws.Columns.Item(8).ColumnWidth = $ws.Columns.Item(5).ColumnWidth

# Set the header text and formulas for the new "CO2/(CO+CO2)" column first
$ws.Range("H1").Value = "CO2/(CO+CO2)"
$ws.Range("H2").Formula = "=F2/(E2+F2)"
$ws.Range("H3:H31").Formula = "=F3/(E3+F3)"

# Now apply the same cell formatting as the neighboring header-style column (I1)
# to every cell in H1:H31 (border, general number format) to match the rest of
# the header/data look used throughout the sheet.
for ($r = 1; $r -le 31; $r++) {
    $ws.Range("I1").Copy()
    $ws.Cells.Item($r, 8).PasteSpecial(-4122)
}

# Restore the active selection to H3, as in the edited workbook
$ws.Range("H3").Select()
